$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.740.86"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.294.41"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'303.52"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").Value = "'96.41"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").Value = "'34.99"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'18.66"
$ws.Range("E12").Value = "  +5.25%  "
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "2.651.98"
$ws.Range("D16").Value = "2.316.22"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "42.654.48"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "'12.80"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "'6.00"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'67.13"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "'236.05"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").Value = "'24.77"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "'167.39"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "'9.00"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "'32.93"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'17.77"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D42").Value = "1.993.74"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "'0.0280"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("D45").Value = "'10.05"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "'2.10"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'2.88"
$ws.Range("E48").Value = "  -5.26%  "
$ws.Range("D49").Value = "'53.62"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").Value = "2.518.85"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "'70.87"
$ws.Range("E51").Value = "  -1.90%  "
